# Generate Report for Archive
# - Update status text from "Ready for handoff" to "In Translation" on all
#   3 sheets (Overview!E2:F2, zh-cn!C2, de-de!C2 all shared the same string).
# - The status column(s) narrow to fit the new, shorter text (Excel
#   auto-sizes these "Status" columns after the content changes).

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$newStatus = "In Translation"

$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsZhCn.Range("C2").Value     = $newStatus
$wsDeDe.Range("C2").Value     = $newStatus

# Narrow the now-shorter "Status" columns to match the new content width.
$newWidth = 12.5

$wsOverview.Columns.Item(5).ColumnWidth = $newWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newWidth
$wsZhCn.Columns.Item(3).ColumnWidth     = $newWidth
$wsDeDe.Columns.Item(3).ColumnWidth     = $newWidth
